$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Test #5 failure/result summary ---
$ws.Range("B7").Value2 = "testFirstConstructor"
$ws.Range("C7").Value2 = "t.getTriangleType()!='S'"
$ws.Range("D7").Value2 = "Expected 'S' but returned 'N'"
$ws.Range("E7").Value2 = "#5"

# --- Row 8: Test #6 failure/result summary (B8 stays blank) ---
$ws.Range("C8").Value2 = "t.getTriangleType()!='S'"
$ws.Range("D8").Value2 = "Expected 'S' but returned 'I'"
$ws.Range("E8").Value2 = "#6"

# --- Row 9: Test #7 failure/result summary (B9 stays blank) ---
$ws.Range("C9").Value2 = "Test compute area"
$ws.Range("D9").Value2 = "Giving unexpected area"
$ws.Range("E9").Value2 = "#7"

# --- Rows 28-29: Fix detail block for bug #5 ---
$ws.Range("B28").Value2 = "#5 Triangle.cpp line 64"
$ws.Range("C28").Value2 = "if (!isTriangle())"
$ws.Range("B29").Value2 = "CHAGNED TO"
$ws.Range("C29").Value2 = "if (isTriangle())"

# --- Rows 31-32: Fix detail block for bug #6 ---
$ws.Range("B31").Value2 = "#6 Triangle.cpp line 76"
$ws.Range("C31").Value2 = "appromixatelyEquals(c,c,m_edgeLengthThreshold))"
$ws.Range("B32").Value2 = "CHANGED TO"
$ws.Range("C32").Value2 = "appromixatelyEquals(c,a,m_edgeLengthThreshold))"

# --- Rows 34-35: Fix detail block for bug #7 ---
$ws.Range("B34").Value2 = "#7 Triangle.cpp line 105"
$ws.Range("C34").Value2 = "double s = ( a + b + b)/2;"
$ws.Range("B35").Value2 = "CHANGED TO"
$ws.Range("C35").Value2 = "double s = ( a + b + c)/2;"

# --- Update view state to match authored workbook ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("C34").Select()
